$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1250.6
$ws.Range("J17").Value = 1273.7838
$ws.Range("L17").Value = 3821.3514
$ws.Range("N17").Value = -4157.3514
$ws.Range("H18").Value = 4207.7856
$ws.Range("I18").Value = 1727.5
$ws.Range("J18").Value = 5199.9
$ws.Range("K18").Value = 1727.5
$ws.Range("L18").Value = 5199.9
$ws.Range("M18").Value = -1443.5
$ws.Range("N18").Value = -5767.9
$ws.Range("H34").Value = 20998.285
$ws.Range("I34").Value = 21998
$ws.Range("K34").Value = 21998
$ws.Range("M34").Value = -21795
$ws.Range("H36").Value = 20998.285
$ws.Range("I36").Value = 21998
$ws.Range("K36").Value = 21998
$ws.Range("M36").Value = -21283
$ws.Range("H38").Value = 2604.2222
$ws.Range("I38").Value = 1408.2
$ws.Range("J38").Value = 4099.25
$ws.Range("K38").Value = 4224.6
$ws.Range("L38").Value = 12297.75
$ws.Range("M38").Value = -3852.6
$ws.Range("N38").Value = -13041.75
$ws.Range("H76").Value = 6497.8
$ws.Range("I76").Value = 4500
$ws.Range("K76").Value = 4500
$ws.Range("M76").Value = -4185
$ws.Range("H79").Value = 6497.8
$ws.Range("I79").Value = 4500
$ws.Range("K79").Value = 4500
$ws.Range("M79").Value = -3408
$ws.Range("H82").Value = 3730.4285
$ws.Range("I82").Value = 3018.8333
$ws.Range("K82").Value = 9056.499899999999
$ws.Range("M82").Value = -8650.499899999999
$ws.Range("H85").Value = 3730.4285
$ws.Range("I85").Value = 3018.8333
$ws.Range("K85").Value = 9056.499899999999
$ws.Range("M85").Value = -7652.499899999999
$ws.Range("H112").Value = 1124.5588
$ws.Range("J112").Value = 1143.0968
$ws.Range("L112").Value = 3429.2904
$ws.Range("N112").Value = -5645.2904
$ws.Range("H132").Value = 38068.426
$ws.Range("I132").Value = 44105.76
$ws.Range("K132").Value = 132317.28
$ws.Range("M132").Value = -129787.28
$ws.Range("H137").Value = 35603644
$ws.Range("J137").Value = 1911280.4
$ws.Range("L137").Value = 5733841.199999999
$ws.Range("N137").Value = -5738941.199999999
$ws.Range("H138").Value = 2253.1
$ws.Range("J138").Value = 2700.7646
$ws.Range("L138").Value = 8102.293799999999
$ws.Range("N138").Value = -18382.2938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6250825
$ws.Range("I32").Value = 6329938
$ws.Range("K32").Value = 6329938
$ws.Range("M32").Value = -6329651
$ws.Range("H45").Value = 2074.3704
$ws.Range("I45").Value = 2045.625
$ws.Range("K45").Value = 2045.625
$ws.Range("M45").Value = -1668.625
$ws.Range("H74").Value = 3907711.5
$ws.Range("I74").Value = 4630547
$ws.Range("K74").Value = 4630547
$ws.Range("M74").Value = -4629673
$ws.Range("H77").Value = 3907711.5
$ws.Range("I77").Value = 4630547
$ws.Range("K77").Value = 23152735
$ws.Range("M77").Value = -23148367
$ws.Range("H97").Value = 1169.3
$ws.Range("I97").Value = 1151.8948
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 1151.8948
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -655.8948
$ws.Range("N97").Value = -2492
$ws.Range("H102").Value = 73348.8
$ws.Range("I102").Value = 90186.25
$ws.Range("K102").Value = 90186.25
$ws.Range("M102").Value = -88564.25
$ws.Range("H135").Value = 90550.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 90550.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 90550.5
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -100690.5
$ws.Range("H138").Value = 90000
$ws.Range("J138").Value = 90000
$ws.Range("L138").Value = 90000
$ws.Range("N138").Value = -100280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1295.25
$ws.Range("I20").Value = 826.75
$ws.Range("J20").Value = 2232.25
$ws.Range("K20").Value = 826.75
$ws.Range("L20").Value = 2232.25
$ws.Range("M20").Value = -579.75
$ws.Range("N20").Value = -2726.25
$ws.Range("H21").Value = 43885.75
$ws.Range("J21").Value = 43885.75
$ws.Range("L21").Value = 43885.75
$ws.Range("N21").Value = -44357.75
$ws.Range("H105").Value = 1416
$ws.Range("I105").Value = 1416
$ws.Range("K105").Value = 1416
$ws.Range("M105").Value = 331

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 92175
$ws.Range("I31").Value = 165495.84
$ws.Range("J31").Value = 5106.5
$ws.Range("K31").Value = 165495.84
$ws.Range("L31").Value = 5106.5
$ws.Range("M31").Value = -165200.84
$ws.Range("N31").Value = -5696.5
$ws.Range("H34").Value = 92175
$ws.Range("I34").Value = 165495.84
$ws.Range("J34").Value = 5106.5
$ws.Range("K34").Value = 165495.84
$ws.Range("L34").Value = 5106.5
$ws.Range("M34").Value = -165293.84
$ws.Range("N34").Value = -5510.5
$ws.Range("H58").Value = 389048.72
$ws.Range("I58").Value = 538072.9399999999
$ws.Range("J58").Value = 8209
$ws.Range("K58").Value = 538072.9399999999
$ws.Range("L58").Value = 8209
$ws.Range("M58").Value = -537869.9399999999
$ws.Range("N58").Value = -8615
$ws.Range("H108").Value = 72224
$ws.Range("I108").Value = 30000
$ws.Range("J108").Value = 80668.8
$ws.Range("K108").Value = 30000
$ws.Range("L108").Value = 80668.8
$ws.Range("M108").Value = -26160
$ws.Range("N108").Value = -88348.8
$ws.Range("H132").Value = 68197110
$ws.Range("I132").Value = 95258264
$ws.Range("K132").Value = 285774792
$ws.Range("M132").Value = -285772262
$ws.Range("H136").Value = 389048.72
$ws.Range("I136").Value = 538072.9399999999
$ws.Range("J136").Value = 8209
$ws.Range("K136").Value = 1614218.82
$ws.Range("L136").Value = 24627
$ws.Range("M136").Value = -1611668.82
$ws.Range("N136").Value = -29727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 500066.66
$ws.Range("J7").Value = 500066.66
$ws.Range("L7").Value = 1500199.98
$ws.Range("N7").Value = -1500423.98
$ws.Range("H74").Value = 9969
$ws.Range("J74").Value = 9969
$ws.Range("L74").Value = 29907
$ws.Range("N74").Value = -32029
$ws.Range("H77").Value = 9969
$ws.Range("J77").Value = 9969
$ws.Range("L77").Value = 89721
$ws.Range("N77").Value = -100329
$ws.Range("H98").Value = 234
$ws.Range("I98").Value = 234
$ws.Range("K98").Value = 702
$ws.Range("M98").Value = 796
$ws.Range("H121").Value = 50001160
$ws.Range("I121").Value = 66666876
$ws.Range("K121").Value = 200000628
$ws.Range("M121").Value = -199999318
$ws.Range("H123").Value = 3875
$ws.Range("I123").Value = 3500
$ws.Range("J123").Value = 5000
$ws.Range("K123").Value = 10500
$ws.Range("L123").Value = 15000
$ws.Range("M123").Value = -8050
$ws.Range("N123").Value = -19900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1715.6471
$ws.Range("I102").Value = 1628.1538
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1628.1538
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -6.153800000000047
$ws.Range("N102").Value = -5244
$ws.Range("H113").Value = 2676.25
$ws.Range("I113").Value = 2438.8333
$ws.Range("J113").Value = 3388.5
$ws.Range("K113").Value = 2438.8333
$ws.Range("L113").Value = 3388.5
$ws.Range("M113").Value = -268.8332999999998
$ws.Range("N113").Value = -7728.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3040.4
$ws.Range("I68").Value = 2400.6667
$ws.Range("K68").Value = 2400.6667
$ws.Range("M68").Value = -1651.6667
$ws.Range("H71").Value = 3040.4
$ws.Range("I71").Value = 2400.6667
$ws.Range("K71").Value = 12003.3335
$ws.Range("M71").Value = -8259.333500000001
$ws.Range("H130").Value = 90429
$ws.Range("J130").Value = 90429
$ws.Range("L130").Value = 90429
$ws.Range("N130").Value = -100469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 4000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H101").Value = 34999.5
$ws.Range("J101").Value = 34999.5
$ws.Range("L101").Value = 34999.5
$ws.Range("N101").Value = -41489.5
